# Weekly refresh of the "Hortaliza, Mapocho Venta Directa de Santiago -
# Zapallo italiano" price series: the data rows (2-13) get re-shuffled as
# a new week's record arrives (rows 6 and 8 are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44277
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 167

$ws.Range("D3").Value = 44186
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 7000
$ws.Range("P3").Value = 117

$ws.Range("D4").Value = 44179
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 7000
$ws.Range("N4").Value = "$/caja 60 unidades"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 117
$ws.Range("Q4").Value = 60

$ws.Range("D5").Value = 44405
$ws.Range("J5").Value = 45
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("N5").Value = "$/caja 50 unidades"
$ws.Range("O5").Value = "Provincia de Quillota"
$ws.Range("P5").Value = 180
$ws.Range("Q5").Value = 50

$ws.Range("D7").Value = 44312
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("P7").Value = 167

$ws.Range("D9").Value = 44284
$ws.Range("J9").Value = 35

$ws.Range("D10").Value = 44315
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("P10").Value = 167

$ws.Range("D11").Value = 44333
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = 10400
$ws.Range("P11").Value = 173

$ws.Range("D12").Value = 44243
$ws.Range("J12").Value = 80
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10375
$ws.Range("O12").Value = "Provincia de Quillota"
$ws.Range("P12").Value = 173

$ws.Range("D13").Value = 44291
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 9000
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = 9000
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 150
